$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF), matching the style of the
# existing header cells (bold / centered / bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-23 for columns I and J.
$values = @(
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(6,7),
    @(8,8),
    @(2,4),
    @(9,9),
    @(8,8),
    @(7,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(1,2),
    @(7,7),
    @(8,8),
    @(5,5),
    @(4,5),
    @(9,9),
    @(6,7),
    @(5,5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
